$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 3107   # F4: 3099 -> 3107
$ws.Cells.Item(6, 6).Value = 258   # F6: 257 -> 258
$ws.Cells.Item(8, 6).Value = 311   # F8: 310 -> 311
$ws.Cells.Item(9, 6).Value = 7295   # F9: 7289 -> 7295
$ws.Cells.Item(10, 6).Value = 63   # F10: 62 -> 63
$ws.Cells.Item(11, 6).Value = 143   # F11: 144 -> 143
$ws.Cells.Item(12, 6).Value = 53   # F12: 50 -> 53
$ws.Cells.Item(14, 6).Value = 23   # F14: 24 -> 23
$ws.Cells.Item(15, 6).Value = 419   # F15: 416 -> 419
$ws.Cells.Item(16, 6).Value = 79   # F16: 78 -> 79
$ws.Cells.Item(17, 6).Value = 79   # F17: 78 -> 79
$ws.Cells.Item(18, 6).Value = 1950   # F18: 1952 -> 1950
$ws.Cells.Item(19, 6).Value = 1775   # F19: 1774 -> 1775
$ws.Cells.Item(20, 6).Value = 1072   # F20: 1071 -> 1072
$ws.Cells.Item(21, 6).Value = 19   # F21: 17 -> 19
$ws.Cells.Item(22, 6).Value = 66   # F22: 65 -> 66
$ws.Cells.Item(23, 6).Value = 1799   # F23: 1791 -> 1799
$ws.Cells.Item(24, 6).Value = 1364   # F24: 1363 -> 1364
$ws.Cells.Item(25, 6).Value = 1221   # F25: 1213 -> 1221
$ws.Cells.Item(26, 6).Value = 632   # F26: 631 -> 632
$ws.Cells.Item(27, 6).Value = 48   # F27: 47 -> 48
$ws.Cells.Item(28, 6).Value = 1113   # F28: 1112 -> 1113
$ws.Cells.Item(29, 6).Value = 115   # F29: 113 -> 115
$ws.Cells.Item(30, 6).Value = 520   # F30: 518 -> 520
$ws.Cells.Item(31, 6).Value = 117   # F31: 110 -> 117
$ws.Cells.Item(32, 6).Value = 65   # F32: 64 -> 65
$ws.Cells.Item(33, 6).Value = 2673   # F33: 2666 -> 2673
$ws.Cells.Item(34, 6).Value = 1505   # F34: 1503 -> 1505
$ws.Cells.Item(35, 6).Value = 2991   # F35: 2984 -> 2991
$ws.Cells.Item(36, 6).Value = 2161   # F36: 2158 -> 2161
$ws.Cells.Item(37, 6).Value = 130   # F37: 129 -> 130
$ws.Cells.Item(41, 6).Value = 37   # F41: 36 -> 37
$ws.Cells.Item(43, 6).Value = 371   # F43: 370 -> 371
$ws.Cells.Item(44, 6).Value = 151   # F44: 148 -> 151
$ws.Cells.Item(45, 6).Value = 508   # F45: 507 -> 508
$ws.Cells.Item(46, 6).Value = 239   # F46: 238 -> 239
$ws.Cells.Item(47, 6).Value = 192   # F47: 190 -> 192
$ws.Cells.Item(48, 6).Value = 725   # F48: 711 -> 725
$ws.Cells.Item(50, 6).Value = 87   # F50: 76 -> 87

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 13   # F4: 12 -> 13
$ws.Cells.Item(13, 6).Value = 55   # F13: 53 -> 55
$ws.Cells.Item(14, 6).Value = 92   # F14: 91 -> 92
$ws.Cells.Item(17, 6).Value = 532   # F17: 531 -> 532
$ws.Cells.Item(20, 6).Value = 28   # F20: 15 -> 28
$ws.Cells.Item(24, 6).Value = 74   # F24: 72 -> 74
$ws.Cells.Item(25, 6).Value = 24   # F25: 23 -> 24
$ws.Cells.Item(31, 6).Value = 7   # F31: 6 -> 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 1831   # F6: 1827 -> 1831
$ws.Cells.Item(7, 6).Value = 1880   # F7: 1879 -> 1880
$ws.Cells.Item(8, 6).Value = 2879   # F8: 2878 -> 2879
$ws.Cells.Item(9, 6).Value = 1121   # F9: 1122 -> 1121
$ws.Cells.Item(10, 6).Value = 1090   # F10: 1088 -> 1090
$ws.Cells.Item(12, 6).Value = 409   # F12: 404 -> 409
$ws.Cells.Item(13, 6).Value = 1808   # F13: 1802 -> 1808
$ws.Cells.Item(14, 6).Value = 8103   # F14: 8090 -> 8103

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 3107   # F3: 3099 -> 3107
$ws.Cells.Item(4, 6).Value = 1831   # F4: 1827 -> 1831
$ws.Cells.Item(5, 6).Value = 311   # F5: 310 -> 311
$ws.Cells.Item(6, 6).Value = 2879   # F6: 2878 -> 2879
$ws.Cells.Item(7, 6).Value = 7295   # F7: 7290 -> 7295
$ws.Cells.Item(8, 6).Value = 1121   # F8: 1122 -> 1121
$ws.Cells.Item(9, 6).Value = 63   # F9: 62 -> 63
$ws.Cells.Item(10, 6).Value = 143   # F10: 144 -> 143
$ws.Cells.Item(11, 6).Value = 409   # F11: 404 -> 409
$ws.Cells.Item(12, 6).Value = 53   # F12: 50 -> 53
$ws.Cells.Item(13, 6).Value = 23   # F13: 24 -> 23
$ws.Cells.Item(14, 6).Value = 420   # F14: 416 -> 420
$ws.Cells.Item(16, 6).Value = 79   # F16: 78 -> 79
$ws.Cells.Item(17, 6).Value = 79   # F17: 78 -> 79
$ws.Cells.Item(19, 6).Value = 1072   # F19: 1071 -> 1072
$ws.Cells.Item(20, 6).Value = 66   # F20: 65 -> 66
$ws.Cells.Item(21, 6).Value = 1799   # F21: 1791 -> 1799
$ws.Cells.Item(22, 6).Value = 1364   # F22: 1363 -> 1364
$ws.Cells.Item(23, 6).Value = 1221   # F23: 1213 -> 1221
$ws.Cells.Item(25, 6).Value = 632   # F25: 631 -> 632
$ws.Cells.Item(26, 6).Value = 48   # F26: 47 -> 48
$ws.Cells.Item(27, 6).Value = 1113   # F27: 1112 -> 1113
$ws.Cells.Item(28, 6).Value = 92   # F28: 91 -> 92
$ws.Cells.Item(29, 6).Value = 532   # F29: 531 -> 532
$ws.Cells.Item(30, 6).Value = 520   # F30: 518 -> 520
$ws.Cells.Item(32, 6).Value = 65   # F32: 64 -> 65
$ws.Cells.Item(33, 6).Value = 2673   # F33: 2666 -> 2673
$ws.Cells.Item(34, 6).Value = 2991   # F34: 2984 -> 2991
$ws.Cells.Item(35, 6).Value = 2161   # F35: 2158 -> 2161
$ws.Cells.Item(40, 6).Value = 37   # F40: 36 -> 37
$ws.Cells.Item(42, 6).Value = 371   # F42: 370 -> 371
$ws.Cells.Item(43, 6).Value = 151   # F43: 148 -> 151
$ws.Cells.Item(44, 6).Value = 74   # F44: 72 -> 74
$ws.Cells.Item(45, 6).Value = 239   # F45: 238 -> 239
